$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# Remove the unused "Sheet1" worksheet (scenarios moved into OrangeLogin).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null

$ws = $wb.Worksheets.Item("OrangeLogin")

# Start from a clean slate on the surviving sheet so no stale
# fill/border formatting from the old table layout carries over.
$ws.Cells.Clear()

# Fourth data column needs a custom width like the other three.
$ws.Columns.Item(4).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "ExpectedError"
$ws.Range("D1").Value = "ErrorType"

$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Font.Size = 12

# ---------------------------------------------------------------------------
# Login test-case rows
# ---------------------------------------------------------------------------
# Row 2: Admin / (blank password) -> Required / empty-pass
$ws.Range("A2").Value = "Admin"
$ws.Range("C2").Value = "Required"
$ws.Range("D2").Value = "empty-pass"

# Row 3: Admin / admin123 -> (no error) / none
$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "admin123"
$ws.Range("D3").Value = "none"

# Row 4: wrong / admin123 -> Invalid credentials / popup
$ws.Range("A4").Value = "wrong"
$ws.Range("B4").Value = "admin123"
$ws.Range("C4").Value = "Invalid credentials"
$ws.Range("D4").Value = "popup"

# Row 5: (blank user) / admin123 -> Required / empty-user
$ws.Range("B5").Value = "admin123"
$ws.Range("C5").Value = "Required"
$ws.Range("D5").Value = "empty-user"

# Row 6: (blank user) / (blank password) -> Required / empty-both
$ws.Range("C6").Value = "Required"
$ws.Range("D6").Value = "empty-both"

# ---------------------------------------------------------------------------
# Restore the active selection on the sheet.
# ---------------------------------------------------------------------------
$ws.Range("K18").Select() | Out-Null
